# Trade #32 closed at 2026-02-17 20:57:28 - unknown UNKNOWN +0.000%
# Applies:
#  - Summary sheet roll-up stats refresh
#  - Strategy Status row for MarketMaking refresh
#  - All Trades: close out trade #60 (row 61) + append new open trade #93 (row 94)
#  - MarketMaking: close out trade #60 (row 28) + append new open trade #93 (row 61)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1400.5
$summary.Range("B4").Value = 0.29
$summary.Range("B5").Value = 0.1
$summary.Range("B6").Value = 60
$summary.Range("B8").Value = 23
$summary.Range("B9").Value = 48.33

# ---------------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 5)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 100.5
$status.Range("D5").Value = 27
$status.Range("E5").Value = 0.18
$status.Range("F5").Value = 0.5
$status.Range("G5").Value = 55.56

# ---------------------------------------------------------------------------
# All Trades sheet - close trade #60 (row 61) and append trade #93 (row 94)
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

$allTrades.Range("G61").Value = 0.2
$allTrades.Range("H61").Value = "CLOSED"
$allTrades.Range("I61").Value = -31.0345
$allTrades.Range("J61").Value = -0.09
$allTrades.Range("K61").Value = 100.5
$allTrades.Range("L61").Value = "early_exit"
$allTrades.Range("M61").Value = 0.14

$allTrades.Range("A94").Value = 93
$allTrades.Range("B94").NumberFormat = "@"
$allTrades.Range("B94").Value = "2026-02-17"
$allTrades.Range("C94").Value = "20:57:21"
$allTrades.Range("D94").Value = "MarketMaking"
$allTrades.Range("E94").Value = "UP"
$allTrades.Range("F94").Value = 0.29
$allTrades.Range("G94").Value = ""
$allTrades.Range("H94").Value = "OPEN"
$allTrades.Range("I94").Value = 0
$allTrades.Range("J94").Value = 0
$allTrades.Range("K94").Value = 100.5855022889912
$allTrades.Range("L94").Value = ""
$allTrades.Range("M94").Value = 0
$allTrades.Range("N94").Value = 0
$allTrades.Range("O94").Value = 0
$allTrades.Range("P94").Value = 0.6
$allTrades.Range("Q94").Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------------
# MarketMaking sheet - close trade #60 (row 28) and append trade #93 (row 61)
# ---------------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")

$mm.Range("G28").Value = 0.2
$mm.Range("H28").Value = "CLOSED"
$mm.Range("I28").Value = -31.0345
$mm.Range("J28").Value = -0.09
$mm.Range("K28").Value = 100.5
$mm.Range("P28").Value = "early_exit"
$mm.Range("Q28").Value = 0.14

$mm.Range("A61").Value = 93
$mm.Range("B61").NumberFormat = "@"
$mm.Range("B61").Value = "2026-02-17"
$mm.Range("C61").Value = "20:57:21"
$mm.Range("D61").Value = "MarketMaking"
$mm.Range("E61").Value = "UP"
$mm.Range("F61").Value = 0.29
$mm.Range("G61").Value = ""
$mm.Range("H61").Value = "OPEN"
$mm.Range("I61").Value = 0
$mm.Range("J61").Value = 0
$mm.Range("K61").Value = 100.5855022889912
$mm.Range("L61").Value = 0
$mm.Range("M61").Value = 0
$mm.Range("N61").Value = 0.6
$mm.Range("O61").Value = "Normal spread capture: 19600 bps"
$mm.Range("P61").Value = ""
$mm.Range("Q61").Value = 0
